$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the product data rows (2-4), keeping the header row and cell styles intact.
$ws.Range("A2:L4").ClearContents()

# Update the selection to reflect the new selected range.
$ws.Range("A2:L4").Select()
